# Swap the match-data (columns B:AD) between pairs of rows.
# Column A (the sequential row id) stays put; everything else in the
# row (match id, teams, scores, odds, ...) is exchanged between the
# two rows in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(112, 114),
    @(116, 117),
    @(118, 120),
    @(121, 122),
    @(155, 156)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B" + $r1 + ":AD" + $r1)
    $rng2 = $ws.Range("B" + $r2 + ":AD" + $r2)

    $vals1 = $rng1.Value()
    $vals2 = $rng2.Value()

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}
